# SASS Learning notebook: add Day 8,9 and Day 10 entries (rows 8 & 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: Day 8, 9 (4-5/07/2019)
$ws.Range("A8").Value = "Day 8, 9 (4-5/07/2019"
$ws.Range("B8").Value = "Learn about new CSS brand feature: background-blend-mode, box-decoration-break"

# New row 9: Day 10 (6/07/2019)
$ws.Range("A9").Value = "Day 10 (6/07/2019"
$ws.Range("B9").Value = "Learn how to use over-flow: hidden when before we used clip-path"

# Match the "Day N" column formatting (centered) used by the rows above
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author left it after typing the new rows
$ws.Range("B13").Select()
